$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.969.73'
$ws.Range('E2').Value = '  -5.40%  '

$ws.Range('D3').Value = '2.920.81'
$ws.Range('E3').Value = '  -3.05%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.76%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '123.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.06%  '

$ws.Range('E7').Value = '  +0.13%  '

$ws.Range('D8').Value = '2.913.36'
$ws.Range('E8').Value = '  -3.24%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.495'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.126'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.25%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '4.80'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -7.33%  '

$ws.Range('E12').Value = '  +2.22%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000214'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.56%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.16'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.34%  '

$ws.Range('E15').Value = '  +1.00%  '

$ws.Range('D16').Value = '3.405.38'
$ws.Range('E16').Value = '  -2.82%  '

$ws.Range('D17').Value = '2.917.59'
$ws.Range('E17').Value = '  -3.00%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.60'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.92%  '

$ws.Range('D19').Value = '57.966.87'
$ws.Range('E19').Value = '  -5.47%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '411.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.54%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.92'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.14%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.668'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.83%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.87'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.80%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.94'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.58%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '77.32'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.11%  '

$ws.Range('E26').Value = '  +0.05%  '

$ws.Range('E27').Value = '  +0.09%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.88%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.31'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.22%  '

$ws.Range('E30').Value = '  +2.59%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.08'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.60%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '24.83'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.74%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0976'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.07%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.916'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.95%  '

$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.41'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.02%  '

$ws.Range('B36').Value = 'Stacks'
$ws.Range('C36').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.01'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -12.08%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '48.16'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.73%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.48'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.26%  '

$ws.Range('D39').Value = '0.0₃0631'
$ws.Range('E39').Value = '  -8.18%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0348'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.62%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.107'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.74%  '

$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.631.39'
$ws.Range('E42').Value = '  -0.72%  '

$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '364.12'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.44%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.41'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.32%  '

$ws.Range('E45').Value = '  -0.02%  '

$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '119.61'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.02%  '

$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.231'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.26%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.97'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.49%  '

$ws.Range('E49').Value = '  +0.57%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.79%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.96'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.61%  '

